$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, matching the styling of the
# other header cells (bold, bordered, centered) by copying E1's format.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the time_taken values for each data row.
$ws.Range("F2").Value = "2021-10-05 13:42:08.684936"
$ws.Range("F3").Value = "2021-10-05 13:42:08.684948"
$ws.Range("F4").Value = "2021-10-05 13:42:08.684952"
$ws.Range("F5").Value = "2021-10-05 13:42:08.684955"
$ws.Range("F6").Value = "2021-10-05 13:42:08.684959"
$ws.Range("F7").Value = "2021-10-05 13:42:08.684962"
$ws.Range("F8").Value = "2021-10-05 13:42:08.684965"
$ws.Range("F9").Value = "2021-10-05 13:42:08.684968"
